# The commit renames the inline-picture "name" attributes (wp:docPr/@name
# and pic:cNvPr/@name) for the three logo pictures living in the document's
# footers/header:
#   - footer1.xml picture (PearsonLogo, id=1/id=0): image2.png -> image1.png
#   - footer2.xml picture (PearsonLogo, id=2/id=0): image2.png -> image1.png
#   - header2.xml picture (BTec_Logo-Orange, id=3/id=0): image1.jpg -> image2.jpg
#
# Word's InlineShape object model has no writable "Name" property (only
# AlternativeText/Title, which map to @descr/@title, are exposed), so the
# rename has to be done by round-tripping the package through
# Document.WordOpenXML and patching the two affected attribute values
# directly in the OOXML text.

$d = $word.ActiveDocument

$xml = $d.WordOpenXML

$xml = $xml.Replace('descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="1" name="image2.png"', 'descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="1" name="image1.png"')
$xml = $xml.Replace('descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image2.png"', 'descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image1.png"')
$xml = $xml.Replace('descr="BTec_Logo-Orange" id="3" name="image1.jpg"', 'descr="BTec_Logo-Orange" id="3" name="image2.jpg"')

# Each of the pictures above repeats the same name on its pic:cNvPr (id="0")
# sibling element immediately afterwards; those use id="0" rather than the
# wp:docPr id, so the three replacements above are combined with these three.
$xml = $xml.Replace('descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image2.png"', 'descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image1.png"')
$xml = $xml.Replace('descr="BTec_Logo-Orange" id="0" name="image1.jpg"', 'descr="BTec_Logo-Orange" id="0" name="image2.jpg"')

$d.WordOpenXML = $xml
